# Modify the NPC Property sheet: mark rows 15-32 of column C ("Public") as TRUE
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column C (Public) to TRUE for rows 15 through 32
for ($row = 15; $row -le 32; $row++) {
    $ws.Cells.Item($row, 3).Value = $true
}

# Update the selection to match the new active range as reflected in the workbook
$ws.Range("C15:C38").Select()
